$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H4").Value = $true
$ws.Range("H5").Value = $true
$ws.Range("G6").Value = $true
$ws.Range("H8").Value = $true
$ws.Range("H15").Value = $true
$ws.Range("H17").Value = $true
$ws.Range("H21").Value = $true
$ws.Range("H29").Value = $true
$ws.Range("H32").Value = $true
$ws.Range("G34").Value = $true
$ws.Range("H36").Value = $true
$ws.Range("H37").Value = $true
$ws.Range("F40").Value = $true
$ws.Range("H42").Value = $true
$ws.Range("H43").Value = $true
$ws.Range("H46").Value = $true
$ws.Range("H48").Value = $true
$ws.Range("H49").Value = $true
$ws.Range("H51").Value = $true
$ws.Range("H55").Value = $true
$ws.Range("B56").Value = "Victor Balaj"
$ws.Range("H56").Value = $true

$ws.Range("M15").Select()
